$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Itga8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.75206033333333
$ws.Range("H2").Value = 38.256181
$ws.Range("I2").Value = 0.1573122343381959
$ws.Range("J2").Value = 0.157312234338196
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3952656666666667
$ws.Range("N2").Value = 1.185797
$ws.Range("O2").Value = 0.02469993039301214
$ws.Range("P2").Value = 0.02469993039301213
$ws.Range("Q2").Value = 5.040451629028555
$ws.Range("R2").Value = 45.364064661257
$ws.Range("S2").Value = 0.003885601238122653
$ws.Range("T2").Value = 0.003885601238122653

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Itga8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.75206033333333
$ws.Range("H3").Value = 38.256181
$ws.Range("I3").Value = 0.1573122343381959
$ws.Range("J3").Value = 0.157312234338196
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.484068333333333
$ws.Range("N3").Value = 10.452205
$ws.Range("O3").Value = 0.2177174811148058
$ws.Range("P3").Value = 0.2177174811148058
$ws.Range("Q3").Value = 44.42904959212277
$ws.Range("R3").Value = 399.861446329105
$ws.Range("S3").Value = 0.03424962340865408
$ws.Range("T3").Value = 0.03424962340865408

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Itga8"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.75206033333333
$ws.Range("H4").Value = 38.256181
$ws.Range("I4").Value = 0.1573122343381959
$ws.Range("J4").Value = 0.157312234338196
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5492906666666667
$ws.Range("N4").Value = 1.647872
$ws.Range("O4").Value = 0.03432486647933305
$ws.Range("P4").Value = 0.03432486647933305
$ws.Range("Q4").Value = 7.004587721870222
$ws.Range("R4").Value = 63.041289496832
$ws.Range("S4").Value = 0.005399721439224127
$ws.Range("T4").Value = 0.005399721439224128

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Itga8"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.75206033333333
$ws.Range("H5").Value = 38.256181
$ws.Range("I5").Value = 0.1573122343381959
$ws.Range("J5").Value = 0.157312234338196
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.574079
$ws.Range("N5").Value = 34.722237
$ws.Range("O5").Value = 0.7232577220128491
$ws.Range("P5").Value = 0.723257722012849
$ws.Range("Q5").Value = 147.5933537107663
$ws.Range("R5").Value = 1328.340183396897
$ws.Range("S5").Value = 0.1137772882521951
$ws.Range("T5").Value = 0.1137772882521951

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Itga8"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.35396833333334
$ws.Range("H6").Value = 61.06190500000001
$ws.Range("I6").Value = 0.2510910513649196
$ws.Range("J6").Value = 0.2510910513649196
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3952656666666667
$ws.Range("N6").Value = 1.185797
$ws.Range("O6").Value = 0.02469993039301214
$ws.Range("P6").Value = 0.02469993039301213
$ws.Range("Q6").Value = 8.045224862587224
$ws.Range("R6").Value = 72.40702376328501
$ws.Range("S6").Value = 0.006201931491021749
$ws.Range("T6").Value = 0.006201931491021748

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Itga8"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 20.35396833333334
$ws.Range("H7").Value = 61.06190500000001
$ws.Range("I7").Value = 0.2510910513649196
$ws.Range("J7").Value = 0.2510910513649196
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.484068333333333
$ws.Range("N7").Value = 10.452205
$ws.Range("O7").Value = 0.2177174811148058
$ws.Range("P7").Value = 0.2177174811148058
$ws.Range("Q7").Value = 70.91461652783612
$ws.Range("R7").Value = 638.231548750525
$ws.Range("S7").Value = 0.05466691123363861
$ws.Range("T7").Value = 0.05466691123363861

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Itga8"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 20.35396833333334
$ws.Range("H8").Value = 61.06190500000001
$ws.Range("I8").Value = 0.2510910513649196
$ws.Range("J8").Value = 0.2510910513649196
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5492906666666667
$ws.Range("N8").Value = 1.647872
$ws.Range("O8").Value = 0.03432486647933305
$ws.Range("P8").Value = 0.03432486647933305
$ws.Range("Q8").Value = 11.18024483512889
$ws.Range("R8").Value = 100.62220351616
$ws.Range("S8").Value = 0.008618666812256221
$ws.Range("T8").Value = 0.008618666812256221

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Itga8"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 20.35396833333334
$ws.Range("H9").Value = 61.06190500000001
$ws.Range("I9").Value = 0.2510910513649196
$ws.Range("J9").Value = 0.2510910513649196
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.574079
$ws.Range("N9").Value = 34.722237
$ws.Range("O9").Value = 0.7232577220128491
$ws.Range("P9").Value = 0.723257722012849
$ws.Range("Q9").Value = 235.5784374534984
$ws.Range("R9").Value = 2120.205937081485
$ws.Range("S9").Value = 0.181603541828003
$ws.Range("T9").Value = 0.181603541828003

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Itga8"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.004706
$ws.Range("H10").Value = 0.014118
$ws.Range("I10").Value = 0.00005805425597465284
$ws.Range("J10").Value = 0.00005805425597465285
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3952656666666667
$ws.Range("N10").Value = 1.185797
$ws.Range("O10").Value = 0.02469993039301214
$ws.Range("P10").Value = 0.02469993039301213
$ws.Range("Q10").Value = 0.001860120227333333
$ws.Range("R10").Value = 0.016741082046
$ws.Range("S10").Value = 0.000001433936081592034
$ws.Range("T10").Value = 0.000001433936081592034

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Vtn"
$ws.Range("C11").Value = "Itga8"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.004706
$ws.Range("H11").Value = 0.014118
$ws.Range("I11").Value = 0.00005805425597465284
$ws.Range("J11").Value = 0.00005805425597465285
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.484068333333333
$ws.Range("N11").Value = 10.452205
$ws.Range("O11").Value = 0.2177174811148058
$ws.Range("P11").Value = 0.2177174811148058
$ws.Range("Q11").Value = 0.01639602557666667
$ws.Range("R11").Value = 0.14756423019
$ws.Range("S11").Value = 0.00001263942637879558
$ws.Range("T11").Value = 0.00001263942637879558

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Vtn"
$ws.Range("C12").Value = "Itga8"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.004706
$ws.Range("H12").Value = 0.014118
$ws.Range("I12").Value = 0.00005805425597465284
$ws.Range("J12").Value = 0.00005805425597465285
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.5492906666666667
$ws.Range("N12").Value = 1.647872
$ws.Range("O12").Value = 0.03432486647933305
$ws.Range("P12").Value = 0.03432486647933305
$ws.Range("Q12").Value = 0.002584961877333334
$ws.Range("R12").Value = 0.023264656896
$ws.Range("S12").Value = 0.000001992704584886982
$ws.Range("T12").Value = 0.000001992704584886982

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Vtn"
$ws.Range("C13").Value = "Itga8"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.004706
$ws.Range("H13").Value = 0.014118
$ws.Range("I13").Value = 0.00005805425597465284
$ws.Range("J13").Value = 0.00005805425597465285
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 11.574079
$ws.Range("N13").Value = 34.722237
$ws.Range("O13").Value = 0.7232577220128491
$ws.Range("P13").Value = 0.723257722012849
$ws.Range("Q13").Value = 0.054467615774
$ws.Range("R13").Value = 0.490208541966
$ws.Range("S13").Value = 0.00004198818892937825
$ws.Range("T13").Value = 0.00004198818892937824

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Vtn"
$ws.Range("C14").Value = "Itga8"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 47.951367
$ws.Range("H14").Value = 143.854101
$ws.Range("I14").Value = 0.5915386600409097
$ws.Range("J14").Value = 0.5915386600409098
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.3952656666666667
$ws.Range("N14").Value = 1.185797
$ws.Range("O14").Value = 0.02469993039301214
$ws.Range("P14").Value = 0.02469993039301213
$ws.Range("Q14").Value = 18.953529044833
$ws.Range("R14").Value = 170.581761403497
$ws.Range("S14").Value = 0.01461096372778614
$ws.Range("T14").Value = 0.01461096372778614

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Vtn"
$ws.Range("C15").Value = "Itga8"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 47.951367
$ws.Range("H15").Value = 143.854101
$ws.Range("I15").Value = 0.5915386600409097
$ws.Range("J15").Value = 0.5915386600409098
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.484068333333333
$ws.Range("N15").Value = 10.452205
$ws.Range("O15").Value = 0.2177174811148058
$ws.Range("P15").Value = 0.2177174811148058
$ws.Range("Q15").Value = 167.065839304745
$ws.Range("R15").Value = 1503.592553742705
$ws.Range("S15").Value = 0.1287883070461343
$ws.Range("T15").Value = 0.1287883070461343

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Vtn"
$ws.Range("C16").Value = "Itga8"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 47.951367
$ws.Range("H16").Value = 143.854101
$ws.Range("I16").Value = 0.5915386600409097
$ws.Range("J16").Value = 0.5915386600409098
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.5492906666666667
$ws.Range("N16").Value = 1.647872
$ws.Range("O16").Value = 0.03432486647933305
$ws.Range("P16").Value = 0.03432486647933305
$ws.Range("Q16").Value = 26.339238347008
$ws.Range("R16").Value = 237.053145123072
$ws.Range("S16").Value = 0.02030448552326781
$ws.Range("T16").Value = 0.02030448552326781

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Vtn"
$ws.Range("C17").Value = "Itga8"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 47.951367
$ws.Range("H17").Value = 143.854101
$ws.Range("I17").Value = 0.5915386600409097
$ws.Range("J17").Value = 0.5915386600409098
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 11.574079
$ws.Range("N17").Value = 34.722237
$ws.Range("O17").Value = 0.7232577220128491
$ws.Range("P17").Value = 0.723257722012849
$ws.Range("Q17").Value = 554.992909815993
$ws.Range("R17").Value = 4994.936188343938
$ws.Range("S17").Value = 0.4278349037437215
$ws.Range("T17").Value = 0.4278349037437215
